$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 4.358119930609447)
    3  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 0, 15.88780690183548)
    4  = @(0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.069511820747223)
    5  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 3.754798637575387)
    6  = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 2.960089034096801)
    7  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    8  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 1, 14.05633640148523)
    9  = @(0.6606524410359556, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 0, 1.345027435901867)
    10 = @(0.2917716402565462, 10.34677158129881, 0.7527432677738641, 10.19245300693656, 0, 21.58373949626578)
    11 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    12 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    13 = @(0.2917716402565462, 0.306821227259698, 22.3905356188092, 1133.036916526867, 0, 1156.026045013193)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
